# Apply edits described by the commit: "logger added, SKU changed, cookies removed"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the existing mistyped value "Shiirtss" -> "Shirts"
$ws.Range("A2").Value = "Shirts"

# Add header row for new columns (SKU, Color, Size)
$ws.Range("B1").Value = "SKU"
$ws.Range("C1").Value = "Color"
$ws.Range("D1").Value = "Size"

# Add data row values (Color, SKU changed to new value, Size)
$ws.Range("C2").Value = "Blue"
$ws.Range("B2").Value = "MSP84FX14025"
$ws.Range("D2").Value = "XL"

# Adjust column widths to match target layout (closest achievable values)
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 15.333333333333334

# Move the active selection to C4 as in the target sheet view
$ws.Range("C4").Select()
